$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$target = "System, dnasr281@gmail.com"
$replacement = "dnasr281@gmail.com, System"

$lastRow = $ws.Cells.SpecialCells(11).Row

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Value2 -eq $target) {
        $cell.Value = $replacement
    }
}
